$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Advance each shop's "remaining days" counter by one day.
# Column E = 剩余 (days remaining), Column F = 开始时间 (cycle start date, yyyyMMdd).
# When E would hit 0 it wraps back to D (total days) and the start date
# rolls forward by D days (a new cycle begins).
for ($row = 2; $row -le 99; $row++) {
    # Row 36 has a malformed start date in the source data and was left untouched.
    if ($row -eq 36) {
        continue
    }

    $remaining = $ws.Cells.Item($row, 5).Value2
    if ($remaining -eq $null) {
        continue
    }

    if ($remaining -eq 1) {
        $total = $ws.Cells.Item($row, 4).Value2
        $startDate = [string]([int]$ws.Cells.Item($row, 6).Value2)
        $year = [int]$startDate.Substring(0, 4)
        $month = [int]$startDate.Substring(4, 2)
        $day = [int]$startDate.Substring(6, 2)
        $newDate = (Get-Date -Year $year -Month $month -Day $day).AddDays($total)

        $ws.Cells.Item($row, 5).Value = $total
        $ws.Cells.Item($row, 6).Value = [int]$newDate.ToString("yyyyMMdd")
    }
    else {
        $ws.Cells.Item($row, 5).Value = $remaining - 1
    }
}
